# "Version extra del word"
#
# Adds, after the first paragraph ("Esta es una prueba mas cañera"):
#   - two empty paragraphs
#   - a new paragraph: "Pues añado una línea mas" (with a gramStart/gramEnd
#     proofErr pair wrapping the word "mas", matching Word's grammar-check
#     markup)
# and, after the existing bookmark paragraph (the one holding the _GoBack
# bookmark) and before the end of the document, one more empty paragraph.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The document currently has exactly two paragraphs:
#   1) "Esta es una prueba mas cañera"
#   2) the (otherwise empty) paragraph carrying the _GoBack bookmark (last
#      paragraph of the document body)
$firstPara = $d.Paragraphs.Item(1)
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Select from the end of paragraph 1 through the end of the bookmark
# paragraph, and replace that whole span with: two empty paragraphs, the
# new "Pues añado..." paragraph, the original bookmark paragraph content
# (preserved as-is), and a trailing empty paragraph.
$targetRange = $d.Range($firstPara.Range.End, $bookmarkPara.Range.End)

$newText = "Pues a" + [char]0x00F1 + "ado una l" + [char]0x00ED + "nea "

$xml = '<w:p ' + $wNs + '/>' `
     + '<w:p ' + $wNs + '/>' `
     + '<w:p ' + $wNs + '>' `
     +   '<w:r><w:t xml:space="preserve">' + $newText + '</w:t></w:r>' `
     +   '<w:proofErr w:type="gramStart"/>' `
     +   '<w:r><w:t>mas</w:t></w:r>' `
     +   '<w:proofErr w:type="gramEnd"/>' `
     + '</w:p>' `
     + '<w:p ' + $wNs + '>' `
     +   '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
     +   '<w:bookmarkEnd w:id="0"/>' `
     + '</w:p>' `
     + '<w:p ' + $wNs + '/>'

$targetRange.InsertXML($xml)
